$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(2)

# Map of shape Id -> new (Left, Top) in points, pre-nudged so that the
# float32 storage used by Shape.Left/Shape.Top round-trips to the exact
# target EMU offsets from the diff (COM Left/Top truncate to EMU on save).
$targets = @{
    7 = @(110.02204895019531, 470.4698486328125)
    8 = @(110.02204895019531, 337.4088439941406)
    9 = @(110.02204895019531, 204.34788513183594)
    10 = @(110.02204895019531, 537.0003662109375)
    11 = @(110.02204895019531, 403.93939208984375)
    12 = @(110.02204895019531, 270.87835693359375)
    13 = @(110.02204895019531, 137.81741333007812)
    16 = @(204.6314239501953, 259.7038879394531)
    17 = @(204.92086791992188, 259.5899353027344)
    18 = @(206.49945068359375, 306.0871887207031)
    19 = @(161.917724609375, 260.944580078125)
    20 = @(186.9940948486328, 295.17144775390625)
    21 = @(204.74954223632812, 291.3163146972656)
    24 = @(337.3480529785156, 470.4698486328125)
    25 = @(337.3480529785156, 337.4088439941406)
    26 = @(337.3480529785156, 204.34788513183594)
    27 = @(337.3480529785156, 537.0003662109375)
    28 = @(337.3480529785156, 403.93939208984375)
    29 = @(337.3480529785156, 270.87835693359375)
    30 = @(337.3480529785156, 137.81741333007812)
    33 = @(507.28143310546875, 472.9762268066406)
    34 = @(521.7424926757812, 446.8822326660156)
    37 = @(564.674072265625, 470.4698486328125)
    38 = @(564.674072265625, 337.4088439941406)
    39 = @(564.674072265625, 204.34788513183594)
    40 = @(564.674072265625, 537.0003662109375)
    41 = @(564.674072265625, 403.93939208984375)
    42 = @(564.674072265625, 270.87835693359375)
    43 = @(564.674072265625, 137.81741333007812)
    46 = @(598.7708740234375, 235.6048126220703)
    47 = @(641.517333984375, 257.4219970703125)
    48 = @(755.6744995117188, 522.3660888671875)
    49 = @(726.7572631835938, 509.2850646972656)
    50 = @(737.8343505859375, 405.89630126953125)
    51 = @(625.1287841796875, 270.7301025390625)
    52 = @(717.0164184570312, 483.1671142578125)
    53 = @(737.3604125976562, 454.6472473144531)
    54 = @(621.9979858398438, 366.3554382324219)
    55 = @(590.6500244140625, 179.0065460205078)
    56 = @(599.4036254882812, 226.0987548828125)
    57 = @(606.3771362304688, 193.92330932617188)
    58 = @(620.8624877929688, 370.9437255859375)
    59 = @(649.8943481445312, 202.82937622070312)
    60 = @(587.302001953125, 175.99961853027344)
    61 = @(619.0615844726562, 175.99757385253906)
    62 = @(592.53564453125, 175.96929931640625)
    63 = @(591.4996948242188, 176.0621337890625)
    64 = @(595.2900390625, 176.0814971923828)
    65 = @(655.661865234375, 385.52764892578125)
    66 = @(653.619873046875, 352.9159851074219)
    67 = @(654.2190551757812, 202.83724975585938)
    68 = @(645.0785522460938, 264.968994140625)
    69 = @(586.8795776367188, 251.88961791992188)
    70 = @(659.2379150390625, 388.35577392578125)
    71 = @(629.8856811523438, 357.2976379394531)
    72 = @(728.3306884765625, 402.1803283691406)
    73 = @(661.7553100585938, 160.39874267578125)
    74 = @(621.0135498046875, 245.96409606933594)
    75 = @(611.6976928710938, 271.6708679199219)
    76 = @(653.9307250976562, 353.8529968261719)
    77 = @(650.4779052734375, 189.37661743164062)
    78 = @(616.7943725585938, 178.9381103515625)
    79 = @(586.787841796875, 281.4774169921875)
    80 = @(639.056640625, 358.95220947265625)
    81 = @(659.6990966796875, 169.8461456298828)
    82 = @(638.0015869140625, 194.34158325195312)
    83 = @(643.7266235351562, 163.9833984375)
    84 = @(599.300048828125, 193.9625244140625)
    85 = @(585.6151733398438, 160.10252380371094)
    86 = @(616.6044311523438, 215.60409545898438)
    87 = @(593.2179565429688, 129.27040100097656)
    88 = @(647.9366455078125, 254.52236938476562)
    89 = @(608.7537231445312, 207.2648162841797)
    90 = @(590.6619873046875, 197.5052032470703)
    91 = @(622.8695068359375, 215.5314178466797)
    92 = @(584.5665893554688, 255.40292358398438)
    93 = @(651.0241088867188, 230.78355407714844)
    94 = @(641.6895751953125, 243.15638732910156)
    95 = @(663.018798828125, 237.4282684326172)
    96 = @(612.1680908203125, 258.7855224609375)
    97 = @(654.420654296875, 296.50189208984375)
    117 = @(90.40811157226562, 533.7175903320312)
    118 = @(90.40811157226562, 400.65655517578125)
    119 = @(90.40811157226562, 267.59552001953125)
    120 = @(90.40811157226562, 134.53457641601562)
    121 = @(107.28228759765625, 537.0003662109375)
    122 = @(107.28228759765625, 403.93939208984375)
    123 = @(107.28228759765625, 270.87835693359375)
    124 = @(107.28228759765625, 137.81741333007812)
}

$applied = 0
for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $sh = $grp.GroupItems.Item($i)
    if ($targets.ContainsKey($sh.Id)) {
        $xy = $targets[$sh.Id]
        $sh.Left = $xy[0]
        $sh.Top = $xy[1]
        $applied++
    }
}

Write-Host "Applied offset updates to" $applied "shapes (expected 89)."
